$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Solved over 400 problems", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start + 12
$digitRng = $d.Range($start, $start + 1)
$digitRng.Text = "5"
$digitRng2 = $d.Range($start, $start + 1)
$digitRng2.Bold = 1
$digitRng3 = $d.Range($start, $start + 1)
$digitRng3.Bold = 0
